$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.516.54"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "'1.872.71"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'315.20"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "'0.5072"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'0.08351"
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").Value = "'41.75"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'6.214"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "'1.872.86"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "'20.38"
$ws.Range("E14").Value = "  +0.18%  "
$ws.Range("D15").Value = "'7.230"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "'1.008"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "'0.00001102"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'91.10"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'0.06699"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "'17.70"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'1.007"
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").Value = "'5.931"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "'28.548.54"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'11.07"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("D26").Value = "'2.084.35"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'161.72"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "'20.62"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'2.358"
$ws.Range("E29").Value = "  -3.10%  "
$ws.Range("D30").Value = "'126.03"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  -2.00%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "'5.789"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  -0.59%  "
$ws.Range("D35").Value = "'0.02450"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "'0.06545"
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("D37").Value = "'0.2159"
$ws.Range("E37").Value = "  -1.08%  "
$ws.Range("D38").Value = "'8.860"
$ws.Range("E38").Value = "  -3.60%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").Value = "'1.251"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'0.6420"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "'0.6030"
$ws.Range("D46").Value = "'13.04"
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("D47").Value = "'3.690"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("D48").Value = "'2.009"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'121.92"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "'1.176"
$ws.Range("E51").Value = "  -8.75%  "
